$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1421.7333
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 1563.5385
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 4690.6155
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -5026.6155

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1872
$ws.Range("I32").Value = 1749.5
$ws.Range("J32").Value = 1994.5
$ws.Range("K32").Value = 1749.5
$ws.Range("L32").Value = 1994.5
$ws.Range("M32").Value = -1423.5
$ws.Range("N32").Value = -2646.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7125
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 7125
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 7125
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -8093

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1200.6316
$ws.Range("J111").Value = 1667.909
$ws.Range("L111").Value = 5003.727000000001
$ws.Range("N111").Value = -11137.727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 14646.1875
$ws.Range("I132").Value = 14646.1875
$ws.Range("K132").Value = 43938.5625
$ws.Range("M132").Value = -41408.5625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4179.231
$ws.Range("I137").Value = 3117.3333
$ws.Range("J137").Value = 5089.4287
$ws.Range("K137").Value = 9351.999899999999
$ws.Range("L137").Value = 15268.2861
$ws.Range("M137").Value = -6801.999899999999
$ws.Range("N137").Value = -20368.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3000
$ws.Range("I26").Value = 3000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 3000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -2670
$ws.Range("N26").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 908.8205
$ws.Range("I32").Value = 576.3333
$ws.Range("K32").Value = 576.3333
$ws.Range("M32").Value = -289.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 13614
$ws.Range("J44").Value = 13614
$ws.Range("L44").Value = 13614
$ws.Range("N44").Value = -14590

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2395
$ws.Range("I99").Value = 2092
$ws.Range("K99").Value = 2092
$ws.Range("M99").Value = -594

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5396.909
$ws.Range("I107").Value = 4052.2856
$ws.Range("K107").Value = 4052.2856
$ws.Range("M107").Value = -2132.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3075.2666
$ws.Range("I58").Value = 2592.1667
$ws.Range("J58").Value = 5007.6665
$ws.Range("K58").Value = 2592.1667
$ws.Range("L58").Value = 5007.6665
$ws.Range("M58").Value = -2389.1667
$ws.Range("N58").Value = -5413.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2319.4375
$ws.Range("I99").Value = 1111
$ws.Range("K99").Value = 1111
$ws.Range("M99").Value = 387

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1800
$ws.Range("J122").Value = 1807.8
$ws.Range("L122").Value = 5423.4
$ws.Range("N122").Value = -10323.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2319.4375
$ws.Range("I126").Value = 1111
$ws.Range("K126").Value = 3333
$ws.Range("M126").Value = -863

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3850.158
$ws.Range("I132").Value = 2297.7144
$ws.Range("J132").Value = 8197
$ws.Range("K132").Value = 6893.1432
$ws.Range("L132").Value = 24591
$ws.Range("M132").Value = -4363.1432
$ws.Range("N132").Value = -29651

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3075.2666
$ws.Range("I136").Value = 2592.1667
$ws.Range("J136").Value = 5007.6665
$ws.Range("K136").Value = 7776.500100000001
$ws.Range("L136").Value = 15022.9995
$ws.Range("M136").Value = -5226.500100000001
$ws.Range("N136").Value = -20122.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 986.6
$ws.Range("J5").Value = 974.8077
$ws.Range("L5").Value = 2924.4231
$ws.Range("N5").Value = -3148.4231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 10
$ws.Range("I29").Value = 10
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 30
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 247
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2666.4167
$ws.Range("J132").Value = 4166.3335
$ws.Range("L132").Value = 37497.0015
$ws.Range("N132").Value = -42557.0015

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 986.6
$ws.Range("J135").Value = 974.8077
$ws.Range("L135").Value = 8773.2693
$ws.Range("N135").Value = -13843.2693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 1681.4
$ws.Range("I22").Value = 453.5
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 453.5
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = 75.5
$ws.Range("N22").Value = -3558

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1921.091
$ws.Range("I122").Value = 1913.2
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5739.6
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3289.6
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 2899.8
$ws.Range("I26").Value = 2124.75
$ws.Range("K26").Value = 2124.75
$ws.Range("M26").Value = -1829.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2714.111
$ws.Range("I40").Value = 2714.111
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2714.111
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2578.111
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2084
$ws.Range("I96").Value = 2502.5
$ws.Range("J96").Value = 1874.75
$ws.Range("K96").Value = 2502.5
$ws.Range("L96").Value = 1874.75
$ws.Range("M96").Value = -1129.5
$ws.Range("N96").Value = -4620.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 855.6923
$ws.Range("J113").Value = 831.7143
$ws.Range("L113").Value = 2495.1429
$ws.Range("N113").Value = -6835.1429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 199497.25
$ws.Range("J141").Value = 99329.664
$ws.Range("L141").Value = 99329.664
$ws.Range("N141").Value = -109689.664
